# Slide 1's notes page (the "Note from Chuck ..." acknowledgement note) is
# translated into Greek. The note is stored in the speaker-notes TextFrame of
# the notes page's body placeholder (Shape 200).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$notesPage = $s.NotesPage
$noteShape = $notesPage.Shapes.Item(1)

$noteText = "Σημείωση από τον  Chuck. Εάν χρησιμοποιείτε αυτό το υλικό, μπορείτε να αφαιρέσετε το λογότυπο UM και να το αντικαταστήσετε με το δικό σας, αλλά διατηρήστε το λογότυπο CC-BY στην πρώτη σελίδα καθώς την/τις σελίδα/ες αναγνώρισης."

$noteShape.TextFrame.TextRange.Text = $noteText
